# Updated cryptos list on Mon Jan 22 21:00:40 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/volume columns (D, E) store plain text (e.g. "308.18", "  -3.75%  ").
# A naive "$cell.Value = '308.18'" lets Excel auto-coerce the numeric-looking
# text into a real number (dropping significant trailing zeros like
# "85.50" -> 85.5, and bumping the cell's style index). Forcing the cell to
# the "@" (Text) number format before the write - then restoring the
# original Style object afterwards - keeps both the literal text and the
# original (unstyled) cell style intact.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "40.135.78"
Set-TextValue $ws.Range("E2") "  -3.61%  "
Set-TextValue $ws.Range("D3") "2.338.17"
Set-TextValue $ws.Range("E3") "  -5.36%  "
Set-TextValue $ws.Range("E4") "  +0.03%  "
Set-TextValue $ws.Range("D5") "308.18"
Set-TextValue $ws.Range("E5") "  -3.75%  "
Set-TextValue $ws.Range("D6") "85.50"
Set-TextValue $ws.Range("E6") "  -6.83%  "
Set-TextValue $ws.Range("D7") "0.529"
Set-TextValue $ws.Range("E7") "  -3.61%  "
Set-TextValue $ws.Range("E8") "  +0.04%  "
Set-TextValue $ws.Range("D9") "0.486"
Set-TextValue $ws.Range("E9") "  -4.58%  "
Set-TextValue $ws.Range("D10") "0.0819"
Set-TextValue $ws.Range("E10") "  -3.83%  "
Set-TextValue $ws.Range("D11") "30.31"
Set-TextValue $ws.Range("E11") "  -7.69%  "
Set-TextValue $ws.Range("E12") "  +0.53%  "
Set-TextValue $ws.Range("D13") "2.701.68"
Set-TextValue $ws.Range("E13") "  -5.30%  "
Set-TextValue $ws.Range("D14") "6.45"
Set-TextValue $ws.Range("E14") "  -5.97%  "
Set-TextValue $ws.Range("D15") "14.75"
Set-TextValue $ws.Range("E15") "  -4.30%  "
Set-TextValue $ws.Range("D16") "2.355.53"
Set-TextValue $ws.Range("E16") "  -4.47%  "
Set-TextValue $ws.Range("D17") "0.757"
Set-TextValue $ws.Range("E17") "  -4.02%  "
Set-TextValue $ws.Range("D18") "40.171.03"
Set-TextValue $ws.Range("E18") "  -3.42%  "
Set-TextValue $ws.Range("D19") "0.0₃0905"
Set-TextValue $ws.Range("E19") "  -3.53%  "
Set-TextValue $ws.Range("D20") "6.10"
Set-TextValue $ws.Range("E20") "  -5.00%  "
Set-TextValue $ws.Range("D21") "67.70"
Set-TextValue $ws.Range("E21") "  -5.05%  "
Set-TextValue $ws.Range("D22") "10.77"
Set-TextValue $ws.Range("E22") "  -3.71%  "
Set-TextValue $ws.Range("D23") "235.65"
Set-TextValue $ws.Range("E23") "  -1.21%  "
Set-TextValue $ws.Range("E24") "  -6.67%  "
Set-TextValue $ws.Range("E25") "  +0.13%  "
Set-TextValue $ws.Range("D26") "1.81"
Set-TextValue $ws.Range("E26") "  -6.66%  "
Set-TextValue $ws.Range("D27") "23.50"
Set-TextValue $ws.Range("E27") "  -5.49%  "
Set-TextValue $ws.Range("D28") "2.21"
Set-TextValue $ws.Range("E28") "  -0.99%  "
Set-TextValue $ws.Range("D29") "9.25"
Set-TextValue $ws.Range("E29") "  -4.88%  "
Set-TextValue $ws.Range("D30") "35.07"
Set-TextValue $ws.Range("E30") "  -3.81%  "
Set-TextValue $ws.Range("D31") "152.30"
Set-TextValue $ws.Range("E31") "  -2.79%  "
Set-TextValue $ws.Range("E32") "  +0.04%  "
Set-TextValue $ws.Range("D33") "5.14"
Set-TextValue $ws.Range("E33") "  -5.05%  "
Set-TextValue $ws.Range("D34") "0.0727"
Set-TextValue $ws.Range("E34") "  -4.86%  "
Set-TextValue $ws.Range("E35") "  -4.55%  "
Set-TextValue $ws.Range("E36") "  -1.76%  "
Set-TextValue $ws.Range("E37") "  -2.60%  "
Set-TextValue $ws.Range("D38") "15.83"
Set-TextValue $ws.Range("E38") "  -6.85%  "
Set-TextValue $ws.Range("D39") "2.75"
Set-TextValue $ws.Range("E39") "  -3.72%  "
Set-TextValue $ws.Range("D40") "1.71"
Set-TextValue $ws.Range("E40") "  -6.35%  "
Set-TextValue $ws.Range("D41") "3.82"
Set-TextValue $ws.Range("E41") "  -4.28%  "
Set-TextValue $ws.Range("D42") "2.28"
Set-TextValue $ws.Range("E42") "  -5.98%  "
Set-TextValue $ws.Range("D43") "1.950.20"
Set-TextValue $ws.Range("E43") "  -2.38%  "
Set-TextValue $ws.Range("D44") "0.0267"
Set-TextValue $ws.Range("E44") "  -5.05%  "
Set-TextValue $ws.Range("D45") "17.76"
Set-TextValue $ws.Range("E45") "  -4.77%  "
Set-TextValue $ws.Range("D46") "9.32"
Set-TextValue $ws.Range("E46") "  -1.07%  "
Set-TextValue $ws.Range("E47") "  -9.02%  "
Set-TextValue $ws.Range("D48") "2.570.40"
Set-TextValue $ws.Range("E48") "  -6.02%  "
Set-TextValue $ws.Range("D49") "92.88"
Set-TextValue $ws.Range("E49") "  -4.51%  "
Set-TextValue $ws.Range("D50") "71.41"
Set-TextValue $ws.Range("E50") "  -5.77%  "
Set-TextValue $ws.Range("D51") "50.72"
Set-TextValue $ws.Range("E51") "  -2.18%  "
